$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L4").Value = 2023
$ws.Range("L5").Value = 33.991563806511245
$ws.Range("L6").Value = 43.352267904134116
$ws.Range("L7").Value = 46.016552065013244
$ws.Range("L8").Value = 57.950845675564537
$ws.Range("L9").Value = 46.481788079470263
$ws.Range("L10").Value = 45.080578284701389
$ws.Range("L11").Value = 39.506289942950417
$ws.Range("L12").Value = 26.964612178240138
$ws.Range("L13").Value = 15.46142526802614
$ws.Range("L14").Value = 33.453947368420813
